# "update to custompromt in TEST"
#
# 1) survey sheet: a new screen-1 question is inserted before the first
#    "end screen" row -> a custom_date prompt named/labelled "test".
#    This pushes every following row down by one (Excel "Insert Row").
# 2) model sheet: a matching new session variable ("test", type "text",
#    isSessionVariable = FALSE) is appended right after the existing ones.
# 3) Selection/active-sheet bookkeeping follows what the author last
#    clicked on: cursor left on F4 of survey, but the workbook's active
#    tab ends up on "model" with B5 selected there.

$wb = $excel.ActiveWorkbook

# --- survey: insert the new custom_date question as row 4 ---
$survey = $wb.Worksheets.Item("survey")
$survey.Rows.Item(4).Insert()
$survey.Cells.Item(4, 4).Value = "custom_date"
$survey.Cells.Item(4, 6).Value = "test"
$survey.Cells.Item(4, 7).Value = "test"

# --- model: register the new "test" session variable ---
$model = $wb.Worksheets.Item("model")
$model.Cells.Item(4, 1).Value = "test"
$model.Cells.Item(4, 2).Value = "text"
$model.Cells.Item(4, 3).Value = $false

# --- selection / active sheet as left by the author ---
[void]$survey.Range("F4").Select()
[void]$model.Select()
[void]$model.Range("B5").Select()
